# The deck ships two DrawingML themes: the theme actually applied to the
# slide master/slides ("Integral" / "Red Violet" palette) and a second,
# unused theme kept only for the notes master ("Office Theme" / "Office"
# palette). The edit swaps which palette is the "live" one: the design
# that drives the slides goes back to the plain Office colour palette.
#
# PowerPoint's automation surface doesn't expose a way to rename a
# ThemeColorScheme or swap the two underlying theme parts wholesale, so we
# reproduce the net visible effect the same way a user would from the
# Design > Variants > Colors > Customize Colors dialog: push each of the
# twelve theme colour slots (Background/Text 1-2, Accent 1-6, Hyperlink,
# Followed Hyperlink) from the Integral values to the Office values, in
# order.

function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # 1  dk1      (Background Dark 1)
    "FFFFFF",  # 2  lt1      (Background Light 1)
    "44546A",  # 3  dk2      (Background Dark 2)
    "E7E6E6",  # 4  lt2      (Background Light 2)
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$colorScheme = $design.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = ToRGB($officeColors[$i - 1])
}
